$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Test 2"
$ws.Range("B2").Value = 20

$ws.Range("A3").Value = "test 333"
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = "test4"
$ws.Range("B4").Value = 0

$ws.Range("A5").Value = "Test567"
$ws.Range("B5").Value = 7
